# #5: cash & deposit done
# Fill in the 存款 (deposits) sheet header row and the new
# property_category/category/date/legislator_name/legislator_id/
# source_file/index columns (G:M) for every existing data row, matching
# the pattern already used on the other property sheets (土地, 建物,
# 汽車, 股票, 保險).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -----------------------------------------------
# Previously row 1 was an accidental duplicate of the first data row;
# give it real column headers instead, and extend through the new
# columns.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

$ws.Range("B1:M1").Font.Bold = $true
$ws.Range("B1:M1").Borders.LineStyle = 1
$ws.Range("B1:M1").HorizontalAlignment = -4108

# --- Data rows (rows 2-8) ----------------------------------------------
# Columns A-F already hold the correct data (index, bank, deposit_type,
# currency, owner, total). Add the metadata columns G-M that every other
# property sheet already carries.
for ($r = 2; $r -le 8; $r++) {
    $idx = $ws.Range("A$r").Value2

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # Force the date column to stay plain text ("2012-03-14") instead of
    # being auto-converted to a date serial number.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2012-03-14"

    $ws.Range("J$r").Value = "林德福"
    $ws.Range("K$r").Value = 908
    $ws.Range("L$r").Value = "tmp82d01"
    $ws.Range("M$r").Value = $idx
}
